$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the stent/balloon detail columns (F:S -> only G:S are extra,
#    F used to be Time_Stent_1). Everything from G through S is dropped;
#    F gets repurposed below as the new "Status" column.
# ---------------------------------------------------------------------------
$ws.Range("G1:S22").EntireColumn.Delete()

# ---------------------------------------------------------------------------
# 2. New "Status" column in F.  F4 used to carry a time-formatted style
#    (it held "Time_Stent_1" data) - reset it to Normal so the new text
#    value isn't dragged into a time/number format.
# ---------------------------------------------------------------------------
$ws.Range("F4").Style = "Normal"
$ws.Cells.Item(1, 6).Value = "Status"

$processedRows = @(2,3,4,5,6,7,8,9,11,12,13,14,15,16,17,18,19,20,21,22)
foreach ($r in $processedRows) {
    $ws.Cells.Item($r, 6).Value = "processed"
}
$ws.Cells.Item(10, 6).Value = "error"

# ---------------------------------------------------------------------------
# 3. Pilot / new patient IDs appended to the log (rows 23-31), all marked
#    "processed" except the last one which is still "new".
# ---------------------------------------------------------------------------
$newIds = @("07027","07032","07036","07037","07045","07052","07053","07063","07071")
$row = 23
foreach ($id in $newIds) {
    $ws.Cells.Item($row, 1).Value = $id
    $ws.Cells.Item($row, 6).Value = "processed"
    $row++
}
# Last added row is still pending processing.
$ws.Cells.Item(31, 6).Value = "new"

# ---------------------------------------------------------------------------
# 4. Freeze the ID column and leave the selection where data entry left off.
# ---------------------------------------------------------------------------
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A32").Select()
